# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 3 (pushing the existing
# rows 3-23 down to 4-24) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3..23 down to 4..24, leaving a blank row 3 behind.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the latest weekly observation.
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 45050
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100104
$ws.Range("H3").Value = "Frutos de pepita"
$ws.Range("I3").Value = 100104003
$ws.Range("J3").Value = "Membrillo"
$ws.Range("K3").Value = "Champion"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 140
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11429
$ws.Range("Q3").Value = "$/caja 18 kilos empedrada"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 635
$ws.Range("T3").Value = 18
